$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 511 (existing rows 511.. shift down to 513..)
$ws.Rows.Item(511).Resize(2).Insert()

$ws.Cells.Item(511, 1).Value = 3
$ws.Cells.Item(511, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(511, 3).Value = "Coquimbo"
$ws.Cells.Item(511, 4).Value = 44637
$ws.Cells.Item(511, 5).Value = 5
$ws.Cells.Item(511, 6).Value = 100114014
$ws.Cells.Item(511, 7).Value = "Betarraga"
$ws.Cells.Item(511, 8).Value = "Sin especificar"
$ws.Cells.Item(511, 9).Value = "Primera"
$ws.Cells.Item(511, 10).Value = 2400
$ws.Cells.Item(511, 11).Value = 550
$ws.Cells.Item(511, 12).Value = 600
$ws.Cells.Item(511, 13).Value = 575
$ws.Cells.Item(511, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(511, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(511, 16).Value = 144
$ws.Cells.Item(511, 17).Value = 4
$ws.Cells.Item(511, 18).Value = "Hortaliza"

$ws.Cells.Item(512, 1).Value = 3
$ws.Cells.Item(512, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(512, 3).Value = "Coquimbo"
$ws.Cells.Item(512, 4).Value = 44637
$ws.Cells.Item(512, 5).Value = 5
$ws.Cells.Item(512, 6).Value = 100114014
$ws.Cells.Item(512, 7).Value = "Betarraga"
$ws.Cells.Item(512, 8).Value = "Sin especificar"
$ws.Cells.Item(512, 9).Value = "Segunda"
$ws.Cells.Item(512, 10).Value = 1200
$ws.Cells.Item(512, 11).Value = 450
$ws.Cells.Item(512, 12).Value = 450
$ws.Cells.Item(512, 13).Value = 450
$ws.Cells.Item(512, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(512, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(512, 16).Value = 112
$ws.Cells.Item(512, 17).Value = 4
$ws.Cells.Item(512, 18).Value = "Hortaliza"
